$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in row 2
$ws.Range("A2").Value = 5001
$ws.Range("B2").Value = 5050

# Fix row 7 values - convert from text "4000.0" to numeric 4000
$ws.Range("A7").Value = 4000
$ws.Range("B7").Value = 4000

# Add new row 8
$ws.Range("A8").Value = 8000
$ws.Range("B8").Value = 8000
